$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force each target cell to Text format before assigning, then clear the
# formatting again so the resulting cell keeps the original (unstyled)
# appearance while still storing the value as a literal string -- this
# avoids Excel auto-converting numeric-looking text (e.g. "0.9887") into a
# real number.
function Set-TextValue($rangeAddr, $text) {
    $rng = $ws.Range($rangeAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
}

Set-TextValue "D2" '27.736.19'
Set-TextValue "E2" '  -0.87%  '
Set-TextValue "D3" '1.849.12'
Set-TextValue "E3" '  -0.47%  '
Set-TextValue "D4" '0.9887'
Set-TextValue "E4" '  -1.49%  '
Set-TextValue "D5" '308.46'
Set-TextValue "E5" '  -1.09%  '
Set-TextValue "D6" '0.9884'
Set-TextValue "E6" '  -1.44%  '
Set-TextValue "D7" '0.5044'
Set-TextValue "E7" '  -0.84%  '
Set-TextValue "D8" '0.3846'
Set-TextValue "E8" '  +1.05%  '
Set-TextValue "D9" '0.08237'
Set-TextValue "E9" '  -0.43%  '
Set-TextValue "D10" '1.103'
Set-TextValue "E10" '  -0.56%  '
Set-TextValue "D11" '40.88'
Set-TextValue "E11" '  -1.73%  '
Set-TextValue "D12" '6.139'
Set-TextValue "E12" '  -1.08%  '
Set-TextValue "B13" 'WrappedEther'
Set-TextValue "C13" 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue "D13" '1.851.72'
Set-TextValue "E13" '  -0.13%  '
Set-TextValue "B14" 'Solana'
Set-TextValue "C14" 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextValue "D14" '20.30'
Set-TextValue "E14" '  -0.95%  '
Set-TextValue "D15" '7.145'
Set-TextValue "E15" '  -0.76%  '
Set-TextValue "D16" '0.9831'
Set-TextValue "E16" '  -2.01%  '
Set-TextValue "D17" '0.00001086'
Set-TextValue "E17" '  -1.02%  '
Set-TextValue "D18" '89.74'
Set-TextValue "E18" '  -0.97%  '
Set-TextValue "E19" '  -0.25%  '
Set-TextValue "D20" '17.49'
Set-TextValue "E20" '  -0.70%  '
Set-TextValue "D21" '0.9875'
Set-TextValue "E21" '  -1.47%  '
Set-TextValue "D22" '5.914'
Set-TextValue "E22" '  -1.77%  '
Set-TextValue "D23" '27.796.58'
Set-TextValue "E23" '  -0.69%  '
Set-TextValue "D24" '10.97'
Set-TextValue "E24" '  -0.55%  '
Set-TextValue "D25" '2.222'
Set-TextValue "E25" '  -0.78%  '
Set-TextValue "D26" '156.88'
Set-TextValue "E26" '  -0.74%  '
Set-TextValue "D27" '2.433'
Set-TextValue "E27" '  -4.08%  '
Set-TextValue "D28" '20.41'
Set-TextValue "E28" '  -0.10%  '
Set-TextValue "D29" '124.59'
Set-TextValue "E29" '  +0.19%  '
Set-TextValue "D30" '0.1045'
Set-TextValue "E30" '  -0.63%  '
Set-TextValue "D31" '1.022'
Set-TextValue "E31" '  -1.16%  '
Set-TextValue "D32" '5.725'
Set-TextValue "E32" '  +1.94%  '
Set-TextValue "D33" '3.548'
Set-TextValue "E33" '  -1.51%  '
Set-TextValue "D34" '9.414'
Set-TextValue "E34" '  +0.03%  '
Set-TextValue "E35" '  +0.09%  '
Set-TextValue "E36" '  -0.37%  '
Set-TextValue "D37" '0.2176'
Set-TextValue "E37" '  +0.61%  '
Set-TextValue "B38" 'TheSandbox'
Set-TextValue "C38" 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue "D38" '0.6390'
Set-TextValue "E38" '  -0.75%  '
Set-TextValue "B39" 'ARBITRUM'
Set-TextValue "C39" 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue "D39" '1.176'
Set-TextValue "E39" '  -2.48%  '
Set-TextValue "D40" '4.907'
Set-TextValue "E40" '  +0.51%  '
Set-TextValue "D41" '1.206'
Set-TextValue "E41" '  -1.49%  '
Set-TextValue "D42" '11.09'
Set-TextValue "E42" '  -0.55%  '
Set-TextValue "D43" '0.6026'
Set-TextValue "E43" '  -0.85%  '
Set-TextValue "D44" '12.91'
Set-TextValue "E44" '  -1.79%  '
Set-TextValue "D45" '1.258'
Set-TextValue "E45" '  -1.66%  '
Set-TextValue "D46" '3.612'
Set-TextValue "E46" '  -1.16%  '
Set-TextValue "D47" '1.981'
Set-TextValue "E47" '  -1.40%  '
Set-TextValue "D48" '1.211'
Set-TextValue "E48" '  +0.47%  '
Set-TextValue "D49" '119.15'
Set-TextValue "E49" '  -0.72%  '
Set-TextValue "E50" '  +0.00%  '
Set-TextValue "D51" '76.96'
Set-TextValue "E51" '  -2.07%  '
